$wb = $excel.ActiveWorkbook

# --- zh-cn sheet (row 6) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("I6").Value = "ffc3f693-f531-4eed-9822-97ef0530bb0e.md"
$wsZh.Range("J6").Value = "ffc3f693-f531-4eed-9822-97ef0530bb0e.46b297c0cf30c60cfdd46604c0c4a0ca39c617d1.zh-cn.xlf"
$wsZh.Range("K6").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6b68f3fa2a379cc2f5e3f23b16827d5ab0f9403c/e2e/ffc3f693-f531-4eed-9822-97ef0530bb0e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b3522a3f9381a31142eb8c37e0faa70a44296b6d/e2e/ffc3f693-f531-4eed-9822-97ef0530bb0e.md."
$wsZh.Range("P6").Value = "2016-10-19 10:47:23"

$wsZh.Hyperlinks.Add($wsZh.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6b68f3fa2a379cc2f5e3f23b16827d5ab0f9403c/e2e/ffc3f693-f531-4eed-9822-97ef0530bb0e.md", [Type]::Missing, [Type]::Missing, "ffc3f693-f531-4eed-9822-97ef0530bb0e.md")

$wsZh.Range("I1:I6").ColumnWidth = 40
$wsZh.Range("J1:J6").ColumnWidth = 40
$wsZh.Range("P1:P6").ColumnWidth = 40

# --- de-de sheet (row 6) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("I6").Value = "ffc3f693-f531-4eed-9822-97ef0530bb0e.md"
$wsDe.Range("J6").Value = "ffc3f693-f531-4eed-9822-97ef0530bb0e.46b297c0cf30c60cfdd46604c0c4a0ca39c617d1.de-de.xlf"
$wsDe.Range("K6").Value = "2016-10-19 10:48:24"
$wsDe.Range("P6").Value = "2016-10-19 10:47:23"

$wsDe.Hyperlinks.Add($wsDe.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6b68f3fa2a379cc2f5e3f23b16827d5ab0f9403c/e2e/ffc3f693-f531-4eed-9822-97ef0530bb0e.md", [Type]::Missing, [Type]::Missing, "ffc3f693-f531-4eed-9822-97ef0530bb0e.md")

$wsDe.Range("I1:I6").ColumnWidth = 40
$wsDe.Range("J1:J6").ColumnWidth = 40
$wsDe.Range("P1:P6").ColumnWidth = 40
